$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 150.77777
$ws.Range("I2").Value = 150.77777
$ws.Range("K2").Value = 150.77777
$ws.Range("M2").Value = -37.77777
$ws.Range("H40").Value = 1356.5938
$ws.Range("I40").Value = 1192.6666
$ws.Range("J40").Value = 1454.95
$ws.Range("K40").Value = 1192.6666
$ws.Range("L40").Value = 1454.95
$ws.Range("M40").Value = -1017.6666
$ws.Range("N40").Value = -1804.95
$ws.Range("H128").Value = 33183.75
$ws.Range("J128").Value = 33183.75
$ws.Range("L128").Value = 33183.75
$ws.Range("N128").Value = -43143.75
$ws.Range("H132").Value = 1667.9642
$ws.Range("I132").Value = 1078.7102
$ws.Range("J132").Value = 4378.533
$ws.Range("K132").Value = 3236.1306
$ws.Range("L132").Value = 13135.599
$ws.Range("M132").Value = -706.1306
$ws.Range("N132").Value = -18195.599
$ws.Range("H137").Value = 2728340.2
$ws.Range("I137").Value = 1283115.8
$ws.Range("J137").Value = 6251075
$ws.Range("K137").Value = 3849347.4
$ws.Range("L137").Value = 18753225
$ws.Range("M137").Value = -3846797.4
$ws.Range("N137").Value = -18758325
$ws.Range("H138").Value = 1348.32
$ws.Range("I138").Value = 457.12195
$ws.Range("J138").Value = 1967.6271
$ws.Range("K138").Value = 1371.36585
$ws.Range("L138").Value = 5902.8813
$ws.Range("M138").Value = 3768.63415
$ws.Range("N138").Value = -16182.8813
$ws.Range("H141").Value = 1265.4366
$ws.Range("I141").Value = 830.3889
$ws.Range("J141").Value = 2647.353
$ws.Range("K141").Value = 2491.1667
$ws.Range("L141").Value = 7942.059
$ws.Range("M141").Value = 2688.8333
$ws.Range("N141").Value = -18302.059
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5208.53
$ws.Range("I32").Value = 3093.8315
$ws.Range("J32").Value = 22318.363
$ws.Range("K32").Value = 3093.8315
$ws.Range("L32").Value = 22318.363
$ws.Range("M32").Value = -2806.8315
$ws.Range("N32").Value = -22892.363
$ws.Range("H45").Value = 3913.9048
$ws.Range("I45").Value = 2733.25
$ws.Range("J45").Value = 5488.1113
$ws.Range("K45").Value = 2733.25
$ws.Range("L45").Value = 5488.1113
$ws.Range("M45").Value = -2356.25
$ws.Range("N45").Value = -6242.1113
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H61").Value = 969
$ws.Range("I61").Value = 961.5682
$ws.Range("J61").Value = 996.25
$ws.Range("K61").Value = 961.5682
$ws.Range("L61").Value = 996.25
$ws.Range("M61").Value = -749.5682
$ws.Range("N61").Value = -1420.25
$ws.Range("H74").Value = 807.0909
$ws.Range("I74").Value = 801.2157
$ws.Range("J74").Value = 882
$ws.Range("K74").Value = 801.2157
$ws.Range("L74").Value = 882
$ws.Range("M74").Value = 72.78430000000003
$ws.Range("N74").Value = -2630
$ws.Range("H77").Value = 807.0909
$ws.Range("I77").Value = 801.2157
$ws.Range("J77").Value = 882
$ws.Range("K77").Value = 4006.0785
$ws.Range("L77").Value = 4410
$ws.Range("M77").Value = 361.9215000000004
$ws.Range("N77").Value = -13146
$ws.Range("H122").Value = 1424.3422
$ws.Range("I122").Value = 1385.1818
$ws.Range("K122").Value = 4155.5454
$ws.Range("M122").Value = -1705.5454
$ws.Range("H136").Value = 969
$ws.Range("I136").Value = 961.5682
$ws.Range("J136").Value = 996.25
$ws.Range("K136").Value = 2884.7046
$ws.Range("L136").Value = 2988.75
$ws.Range("M136").Value = -334.7046
$ws.Range("N136").Value = -8088.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 26000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 26000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 26000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -27872
$ws.Range("H78").Value = 26000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 26000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 78000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -87360
$ws.Range("H141").Value = 22770.715
$ws.Range("I141").Value = 20000
$ws.Range("J141").Value = 23232.5
$ws.Range("K141").Value = 20000
$ws.Range("L141").Value = 23232.5
$ws.Range("M141").Value = -14820
$ws.Range("N141").Value = -33592.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3269.1428
$ws.Range("I31").Value = 2447.5
$ws.Range("J31").Value = 4364.6665
$ws.Range("K31").Value = 2447.5
$ws.Range("L31").Value = 4364.6665
$ws.Range("M31").Value = -2152.5
$ws.Range("N31").Value = -4954.6665
$ws.Range("H34").Value = 3269.1428
$ws.Range("I34").Value = 2447.5
$ws.Range("J34").Value = 4364.6665
$ws.Range("K34").Value = 2447.5
$ws.Range("L34").Value = 4364.6665
$ws.Range("M34").Value = -2245.5
$ws.Range("N34").Value = -4768.6665
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H132").Value = 1526.8861
$ws.Range("I132").Value = 1407.75
$ws.Range("J132").Value = 2035.2
$ws.Range("K132").Value = 4223.25
$ws.Range("L132").Value = 6105.6
$ws.Range("M132").Value = -1693.25
$ws.Range("N132").Value = -11165.6
$ws.Range("H134").Value = 2987.476
$ws.Range("I134").Value = 3115.1455
$ws.Range("J134").Value = 2109.75
$ws.Range("K134").Value = 9345.4365
$ws.Range("L134").Value = 6329.25
$ws.Range("M134").Value = -6810.4365
$ws.Range("N134").Value = -11399.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1971557.2
$ws.Range("I2").Value = 125178.875
$ws.Range("J2").Value = 4081703.8
$ws.Range("K2").Value = 751073.25
$ws.Range("L2").Value = 24490222.8
$ws.Range("M2").Value = -750960.25
$ws.Range("N2").Value = -24490448.8
$ws.Range("H38").Value = 83.59999999999999
$ws.Range("I38").Value = 74
$ws.Range("J38").Value = 94.57143000000001
$ws.Range("K38").Value = 222
$ws.Range("L38").Value = 283.71429
$ws.Range("M38").Value = 125
$ws.Range("N38").Value = -977.71429
$ws.Range("H131").Value = 911.7625
$ws.Range("I131").Value = 516.9
$ws.Range("J131").Value = 968.17145
$ws.Range("K131").Value = 1550.7
$ws.Range("L131").Value = 2904.51435
$ws.Range("M131").Value = 3489.3
$ws.Range("N131").Value = -12984.51435
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 4333.3335
$ws.Range("I44").Value = 4000
$ws.Range("J44").Value = 4500
$ws.Range("K44").Value = 4000
$ws.Range("L44").Value = 4500
$ws.Range("M44").Value = -3404
$ws.Range("N44").Value = -5692
$ws.Range("H80").Value = 1154670
$ws.Range("I80").Value = 2461.8333
$ws.Range("J80").Value = 2750035.2
$ws.Range("K80").Value = 2461.8333
$ws.Range("L80").Value = 2750035.2
$ws.Range("M80").Value = -1463.8333
$ws.Range("N80").Value = -2752031.2
$ws.Range("H83").Value = 1154670
$ws.Range("I83").Value = 2461.8333
$ws.Range("J83").Value = 2750035.2
$ws.Range("K83").Value = 12309.1665
$ws.Range("L83").Value = 13750176
$ws.Range("M83").Value = -7317.166499999999
$ws.Range("N83").Value = -13760160
$ws.Range("H102").Value = 2044.9333
$ws.Range("I102").Value = 1762.4286
$ws.Range("J102").Value = 6000
$ws.Range("K102").Value = 1762.4286
$ws.Range("L102").Value = 6000
$ws.Range("M102").Value = -140.4286
$ws.Range("N102").Value = -9244
$ws.Range("H122").Value = 1507.3
$ws.Range("I122").Value = 1548.1052
$ws.Range("J122").Value = 1436.8182
$ws.Range("K122").Value = 4644.3156
$ws.Range("L122").Value = 4310.4546
$ws.Range("M122").Value = -2194.3156
$ws.Range("N122").Value = -9210.454600000001
$ws.Range("H132").Value = 1549.3925
$ws.Range("I132").Value = 1202.2075
$ws.Range("J132").Value = 2257.1155
$ws.Range("K132").Value = 3606.6225
$ws.Range("L132").Value = 6771.3465
$ws.Range("M132").Value = -1076.6225
$ws.Range("N132").Value = -11831.3465
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 368.3846
$ws.Range("I16").Value = 348.8
$ws.Range("J16").Value = 433.66666
$ws.Range("K16").Value = 348.8
$ws.Range("L16").Value = 433.66666
$ws.Range("M16").Value = -178.8
$ws.Range("N16").Value = -773.66666
$ws.Range("H122").Value = 2798.963
$ws.Range("I122").Value = 2383.6667
$ws.Range("J122").Value = 2917.6191
$ws.Range("K122").Value = 7151.000100000001
$ws.Range("L122").Value = 8752.8573
$ws.Range("M122").Value = -4701.000100000001
$ws.Range("N122").Value = -13652.8573
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 21000
$ws.Range("J70").Value = 21000
$ws.Range("L70").Value = 21000
$ws.Range("N70").Value = -21630
$ws.Range("H73").Value = 21000
$ws.Range("J73").Value = 21000
$ws.Range("L73").Value = 21000
$ws.Range("N73").Value = -23184
$ws.Range("H132").Value = 1192.3658
$ws.Range("I132").Value = 1081.863
$ws.Range("J132").Value = 2088.6667
$ws.Range("K132").Value = 3245.589
$ws.Range("L132").Value = 6266.000100000001
$ws.Range("M132").Value = -715.5889999999999
$ws.Range("N132").Value = -11326.0001
$ws.Range("H136").Value = 1213.9697
$ws.Range("I136").Value = 1213.9697
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3641.9091
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1091.9091
$ws.Range("N136").ClearContents()
